$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3617.8928
$ws.Range("I64").Value = 3277.889
$ws.Range("J64").Value = 4229.9
$ws.Range("K64").Value = 3277.889
$ws.Range("L64").Value = 4229.9
$ws.Range("M64").Value = -3029.889

$ws.Range("H67").Value = 3617.8928
$ws.Range("I67").Value = 3277.889
$ws.Range("J67").Value = 4229.9
$ws.Range("K67").Value = 3277.889
$ws.Range("L67").Value = 4229.9
$ws.Range("M67").Value = -2419.889

$ws.Range("H96").Value = 25002156
$ws.Range("I96").Value = 62502710
$ws.Range("J96").Value = 1784.5
$ws.Range("K96").Value = 187508130
$ws.Range("L96").Value = 5353.5
$ws.Range("M96").Value = -187506757
$ws.Range("N96").Value = -8099.5

$ws.Range("H101").Value = 41667020
$ws.Range("I101").Value = 41667020
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 125001060
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -124999438

$ws.Range("H132").Value = 697.25
$ws.Range("I132").Value = 697.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2091.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 438.25
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 544.8511
$ws.Range("I135").Value = 531.2889
$ws.Range("J135").Value = 850
$ws.Range("K135").Value = 4781.6001
$ws.Range("L135").Value = 7650
$ws.Range("M135").Value = -2246.6001
$ws.Range("N135").Value = -12720

$ws.Range("H138").Value = 1402.1897
$ws.Range("I138").Value = 1151.091
$ws.Range("J138").Value = 6005.6665
$ws.Range("K138").Value = 3453.273
$ws.Range("L138").Value = 18016.9995
$ws.Range("M138").Value = 1686.727
$ws.Range("N138").Value = -28296.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1198
$ws.Range("I61").Value = 997.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 997.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -785.5
$ws.Range("N61").Value = -2424

$ws.Range("H110").Value = 17785.533
$ws.Range("I110").Value = 21755.959
$ws.Range("J110").Value = 1903.8334
$ws.Range("K110").Value = 21755.959
$ws.Range("L110").Value = 1903.8334
$ws.Range("M110").Value = -19710.959
$ws.Range("N110").Value = -5993.8334

$ws.Range("H132").Value = 2002.9111
$ws.Range("I132").Value = 1178.7188
$ws.Range("J132").Value = 4031.6924
$ws.Range("K132").Value = 3536.1564
$ws.Range("L132").Value = 12095.0772
$ws.Range("M132").Value = -1006.1564
$ws.Range("N132").Value = -17155.0772

$ws.Range("H133").Value = 34925.668
$ws.Range("I133").Value = 10000
$ws.Range("J133").Value = 47388.5
$ws.Range("K133").Value = 10000
$ws.Range("L133").Value = 47388.5
$ws.Range("M133").Value = -7470
$ws.Range("N133").Value = -52448.5

$ws.Range("H134").Value = 31651.143
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 31651.143
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 31651.143
$ws.Range("N134").Value = -41791.143

$ws.Range("H135").Value = 83143
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 83143
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 83143
$ws.Range("N135").Value = -93283

$ws.Range("H136").Value = 1198
$ws.Range("I136").Value = 997.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2992.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -442.5
$ws.Range("N136").Value = -11100

$ws.Range("H137").Value = 49970
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 49970
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 49970
$ws.Range("N137").Value = -60170

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1200
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1200
$ws.Range("N107").Value = -5040
$ws.Range("M107").ClearContents()

$ws.Range("H134").Value = 1813.0834
$ws.Range("I134").Value = 1560.3334
$ws.Range("J134").Value = 2571.3333
$ws.Range("K134").Value = 4681.0002
$ws.Range("L134").Value = 7713.999899999999
$ws.Range("M134").Value = -2146.0002
$ws.Range("N134").Value = -12783.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H58").Value = 1300.561
$ws.Range("I58").Value = 1170.2667
$ws.Range("J58").Value = 1655.909
$ws.Range("K58").Value = 1170.2667
$ws.Range("L58").Value = 1655.909
$ws.Range("M58").Value = -967.2666999999999
$ws.Range("N58").Value = -2061.909

$ws.Range("H62").Value = 5601.1
$ws.Range("I62").Value = 4167.5
$ws.Range("J62").Value = 7751.5
$ws.Range("K62").Value = 4167.5
$ws.Range("L62").Value = 7751.5
$ws.Range("M62").Value = -3543.5
$ws.Range("N62").Value = -8999.5

$ws.Range("H65").Value = 5601.1
$ws.Range("I65").Value = 4167.5
$ws.Range("J65").Value = 7751.5
$ws.Range("K65").Value = 20837.5
$ws.Range("L65").Value = 38757.5
$ws.Range("M65").Value = -17717.5
$ws.Range("N65").Value = -44997.5

$ws.Range("H134").Value = 1470.375
$ws.Range("I134").Value = 1603.3823
$ws.Range("J134").Value = 716.6667
$ws.Range("K134").Value = 4810.1469
$ws.Range("L134").Value = 2150.0001
$ws.Range("M134").Value = -2275.1469
$ws.Range("N134").Value = -7220.0001

$ws.Range("H135").Value = 42836.93
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 42836.93
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 42836.93
$ws.Range("N135").Value = -52976.93

$ws.Range("H136").Value = 1300.561
$ws.Range("I136").Value = 1170.2667
$ws.Range("J136").Value = 1655.909
$ws.Range("K136").Value = 3510.800099999999
$ws.Range("L136").Value = 4967.727000000001
$ws.Range("M136").Value = -960.8000999999995
$ws.Range("N136").Value = -10067.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 772.5
$ws.Range("I46").Value = 150
$ws.Range("J46").Value = 980
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 2940
$ws.Range("M46").Value = -359
$ws.Range("N46").Value = -3122

$ws.Range("H131").Value = 4657.6206
$ws.Range("I131").Value = 548.2857
$ws.Range("J131").Value = 5965.136
$ws.Range("K131").Value = 1644.8571
$ws.Range("L131").Value = 17895.408
$ws.Range("M131").Value = 3395.1429
$ws.Range("N131").Value = -27975.408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2769.348
$ws.Range("I132").Value = 2187
$ws.Range("J132").Value = 3861.25
$ws.Range("K132").Value = 6561
$ws.Range("L132").Value = 11583.75
$ws.Range("M132").Value = -4031
$ws.Range("N132").Value = -16643.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H123").Value = 1521500
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 1521500
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 1521500
$ws.Range("N123").Value = -1531300

$ws.Range("H125").Value = 27715
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 27715
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 27715
$ws.Range("N125").Value = -37555

$ws.Range("H127").Value = 72335.336
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 72335.336
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 72335.336
$ws.Range("N127").Value = -82255.336
$ws.Range("M127").ClearContents()

$ws.Range("H130").Value = 48500
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 48500
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 48500
$ws.Range("N130").Value = -58540

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 39199.375
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 39199.375
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 39199.375
$ws.Range("N128").Value = -49159.375

